# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = "'58.043.20"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.52%  '

$ws.Range('D3').Value = "'2.577.88"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.60%  '

$ws.Range('D5').Value = "'518.26"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.62%  '

$ws.Range('D6').Value = "'142.37"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.52%  '

$ws.Range('D7').Value = "'0.997"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.21%  '

$ws.Range('E8').Value = '  -0.76%  '

$ws.Range('D9').Value = "'2.595.63"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.17%  '

$ws.Range('D10').Value = "'6.80"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.55%  '

$ws.Range('E11').Value = '  -1.41%  '

$ws.Range('E12').Value = '  -3.79%  '

$ws.Range('E13').Value = '  -1.16%  '

$ws.Range('D14').Value = "'3.034.19"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.57%  '

$ws.Range('D15').Value = "'58.059.73"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.46%  '

$ws.Range('D16').Value = "'20.35"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.41%  '

$ws.Range('E17').Value = '  -1.82%  '

$ws.Range('D18').Value = "'2.569.73"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.02%  '

$ws.Range('D19').Value = "'343.48"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.40%  '

$ws.Range('D20').Value = "'4.31"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.40%  '

$ws.Range('D21').Value = "'10.29"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.65%  '

$ws.Range('D22').Value = "'6.36"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.08%  '

$ws.Range('E23').Value = '  -0.11%  '

$ws.Range('D24').Value = "'65.75"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.60%  '

$ws.Range('E25').Value = '  -1.62%  '

$ws.Range('D26').Value = "'0.402"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.24%  '

$ws.Range('E27').Value = '  -0.24%  '

$ws.Range('D28').Value = "'2.689.75"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.69%  '

$ws.Range('D29').Value = "'7.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.57%  '

$ws.Range('D30').Value = "'0.0₃0745"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.63%  '

$ws.Range('D32').Value = "'6.24"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.23%  '

$ws.Range('E33').Value = '  -0.40%  '

$ws.Range('D34').Value = "'18.76"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.54%  '

$ws.Range('D35').Value = "'149.73"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.85%  '

$ws.Range('D36').Value = "'4.02"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.93%  '

$ws.Range('E37').Value = '  -3.10%  '

$ws.Range('D38').Value = "'0.868"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.97%  '

$ws.Range('D39').Value = "'36.02"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.19%  '

$ws.Range('D40').Value = "'0.834"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.88%  '

$ws.Range('D41').Value = "'1.45"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.56%  '

$ws.Range('E42').Value = '  -3.16%  '

$ws.Range('D43').Value = "'0.997"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.27%  '

$ws.Range('D44').Value = "'271.86"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.61%  '

$ws.Range('D45').Value = "'10.65"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.19%  '

$ws.Range('D46').Value = "'0.0951"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.80%  '

$ws.Range('D47').Value = "'0.588"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.11%  '

$ws.Range('D48').Value = "'18.86"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.12%  '

$ws.Range('D49').Value = "'0.0524"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.49%  '

# Row 50 - RenderToken inserted in place of Maker (new coin enters the ranking)
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'4.66"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.09%  '

# Row 51 - Maker shifts down, replacing InjectiveProtocol which drops off the list
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = "'1.972.44"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.38%  '
